# Apply the latest cryptocurrency market snapshot to the "cryptos" sheet.
# For each affected row, the Coin/Link/Price/Volume(1h) cells are refreshed
# with freshly scraped values (Price column is forced to Text so that
# values such as "1.10" or "0.0000205" are not reinterpreted as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = "93.212.10"; Volume = "  +1.78%  " },
    @{ Row = 3; Price = "3.124.36"; Volume = "  +0.33%  " },
    @{ Row = 4; Volume = "  +0.06%  " },
    @{ Row = 5; Price = "238.51"; Volume = "  -2.54%  " },
    @{ Row = 6; Price = "614.66"; Volume = "  -0.40%  " },
    @{ Row = 7; Price = "1.10"; Volume = "  +0.31%  " },
    @{ Row = 8; Price = "0.394"; Volume = "  +2.49%  " },
    @{ Row = 9; Volume = "  -0.03%  " },
    @{ Row = 10; Price = "0.813"; Volume = "  +7.81%  " },
    @{ Row = 11; Price = "3.119.92"; Volume = "  +0.19%  " },
    @{ Row = 12; Volume = "  -2.09%  " },
    @{ Row = 13; Volume = "  -1.87%  " },
    @{ Row = 14; Coin = "Avalanche"; Link = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; Price = "34.73"; Volume = "  -0.12%  " },
    @{ Row = 15; Coin = "WrappedBTC"; Link = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; Price = "92.827.93"; Volume = "  +1.37%  " },
    @{ Row = 16; Price = "5.46"; Volume = "  -2.73%  " },
    @{ Row = 17; Price = "3.707.76"; Volume = "  +0.18%  " },
    @{ Row = 18; Price = "3.108.68"; Volume = "  -0.33%  " },
    @{ Row = 19; Volume = "  +1.10%  " },
    @{ Row = 20; Price = "14.63"; Volume = "  -1.51%  " },
    @{ Row = 21; Price = "5.94"; Volume = "  +2.42%  " },
    @{ Row = 22; Price = "0.0000205"; Volume = "  +1.53%  " },
    @{ Row = 23; Price = "442.79"; Volume = "  -1.46%  " },
    @{ Row = 24; Price = "9.18"; Volume = "  -0.15%  " },
    @{ Row = 25; Volume = "  +5.02%  " },
    @{ Row = 26; Price = "5.64"; Volume = "  -3.74%  " },
    @{ Row = 27; Price = "12.32"; Volume = "  +5.39%  " },
    @{ Row = 28; Price = "85.85"; Volume = "  -3.95%  " },
    @{ Row = 29; Volume = "  -0.01%  " },
    @{ Row = 30; Price = "0.180"; Volume = "  +8.39%  " },
    @{ Row = 31; Price = "0.234"; Volume = "  +3.56%  " },
    @{ Row = 32; Price = "0.124"; Volume = "  -9.25%  " },
    @{ Row = 33; Price = "9.21"; Volume = "  -0.88%  " },
    @{ Row = 34; Volume = "  -4.94%  " },
    @{ Row = 35; Price = "8.10"; Volume = "  +5.96%  " },
    @{ Row = 36; Volume = "  -8.92%  " },
    @{ Row = 37; Price = "25.88"; Volume = "  -1.37%  " },
    @{ Row = 38; Price = "3.97"; Volume = "  +1.04%  " },
    @{ Row = 39; Price = "1.91"; Volume = "  -3.46%  " },
    @{ Row = 40; Volume = "  -0.74%  " },
    @{ Row = 41; Price = "23.98"; Volume = "  +8.00%  " },
    @{ Row = 42; Price = "474.56"; Volume = "  -2.75%  " },
    @{ Row = 43; Price = "0.432"; Volume = "  -0.62%  " },
    @{ Row = 44; Volume = "  -2.31%  " },
    @{ Row = 46; Price = "159.76"; Volume = "  -0.23%  " },
    @{ Row = 47; Price = "0.690"; Volume = "  -1.42%  " },
    @{ Row = 48; Price = "1.85"; Volume = "  -2.97%  " },
    @{ Row = 49; Price = "1.33"; Volume = "  -0.98%  " },
    @{ Row = 50; Price = "44.15"; Volume = "  +0.06%  " },
    @{ Row = 51; Coin = "Filecoin"; Link = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Price = "4.38"; Volume = "  -0.40%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("Coin"))   { $ws.Range("B$($u.Row)").Value2 = $u.Coin }
    if ($u.ContainsKey("Link"))   { $ws.Range("C$($u.Row)").Value2 = $u.Link }
    if ($u.ContainsKey("Price")) {
        $priceCell = $ws.Range("D$($u.Row)")
        # Keep the Price column as plain text so values like "1.10",
        # "0.0000205" or trailing zeros survive the round-trip untouched.
        $priceCell.NumberFormat = "@"
        $priceCell.Value2 = $u.Price
    }
    if ($u.ContainsKey("Volume")) { $ws.Range("E$($u.Row)").Value2 = $u.Volume }
}
